# Rename the existing sheet and add a new one after it, populate the new
# rows, and set selection/active tab to match the target state.

$wb = $excel.ActiveWorkbook

# Rename Sheet1 -> "To fix"
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "To fix"

# Add new row 7 with the bug text
$ws1.Range("A7").Value = "remove the requirement of four wrong answers. One should be enough to save the exercise."

# Add a new worksheet right after "To fix", named "to Implement"
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "to Implement"
$ws2.Range("A1").Value = "assiginments: to allow a teacher to create new assignments"

# Set selections to match target state
$ws1.Range("A8").Select()
$ws2.Range("C8").Select()

# Make "to Implement" the active sheet/tab
$ws2.Activate()
